# Fixed Stimulus Absolute Timestamps
# Rename the task-order sheets (new timestamp suffixes) and update the
# per-sheet stimulus-file lists in column B to their corrected values.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778612935307"
$ws1.Range("B2").Value = "go_stims-16504778612545316.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778612765307.csv"
$ws1.Range("B4").Value = "go_stims-16504778612775278.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477861292562.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778640326142"
$ws2.Range("B2").Value = "OB-16504778618525317.csv"
$ws2.Range("B3").Value = "TB-16504778638945618.csv"
$ws2.Range("B4").Value = "TB-16504778640185611.csv"
$ws2.Range("B5").Value = "ZB-match_5-1650477861353565.csv"
$ws2.Range("B6").Value = "ZB-match_8-16504778615475287.csv"
$ws2.Range("B7").Value = "OB-16504778632675278.csv"
$ws2.Range("B8").Value = "ZB-match_7-16504778618175313.csv"
$ws2.Range("B9").Value = "TB-1650477863618559.csv"
$ws2.Range("B10").Value = "OB-16504778623415656.csv"

# --- Sheet 3: RS (name only, no data changes) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650477864038527"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778640975273"
$ws4.Range("B2").Value = "MM_stims-16504778640645602.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778640395272.csv"
$ws4.Range("B4").Value = "MM_stims-16504778640805604.csv"
$ws4.Range("B5").Value = "ZM_stims-165047786406553.csv"
$ws4.Range("B6").Value = "MM_stims-16504778640965595.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778640805604.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778641605313"
$ws5.Range("B2").Value = "vSAT_stims-16504778641445618.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778641285286.csv"
$ws5.Range("B4").Value = "SAT_stims-1650477864112559.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778641005285.csv"
